$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# --- Cell F2: update the rpc-reply message-id UUID ---
$cellF2 = $ws.Range("F2")
$f2 = $cellF2.Value2
$f2 = $f2.Replace("b985be53-4e1a-4d23-acc9-dea3bd1af9b9", "f583f20c-df8c-491c-9e9f-91f73da4bbe5")
$cellF2.Value2 = $f2

# --- Cell G2: update identifier/name elements for the BGP protocol block ---
$cellG2 = $ws.Range("G2")
$g2 = $cellG2.Value2
$g2 = $g2.Replace('<identifier>BGP</identifier>' + "`n" + '              <name>BGP_65000</name>', '<identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:BGP</identifier>' + "`n" + '              <name>default</name>')
$g2 = $g2.Replace('<name>BGP_65000</name>' + "`n" + '              </config>', '<name>default</name>' + "`n" + '              </config>')
$cellG2.Value2 = $g2
